$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 2).Value = 0.09375
$ws.Cells.Item(10, 2).Value = 0.09375
$ws.Cells.Item(11, 2).Value = 0.109375
$ws.Cells.Item(12, 2).Value = 0.109375
$ws.Cells.Item(13, 2).Value = 0.109375
$ws.Cells.Item(14, 2).Value = 0.109375
$ws.Cells.Item(16, 2).Value = 0.234375
$ws.Cells.Item(17, 2).Value = 0.234375
$ws.Cells.Item(18, 2).Value = 0.15625
$ws.Cells.Item(19, 2).Value = 0.140625
$ws.Cells.Item(20, 2).Value = 0.21875
$ws.Cells.Item(21, 2).Value = 0.125
$ws.Cells.Item(22, 2).Value = 0.125
$ws.Cells.Item(23, 2).Value = 0.140625
$ws.Cells.Item(24, 2).Value = 0.140625
$ws.Cells.Item(25, 2).Value = 0.125
$ws.Cells.Item(26, 2).Value = 0.125
$ws.Cells.Item(27, 2).Value = 0.125
$ws.Cells.Item(28, 2).Value = 0.125
$ws.Cells.Item(29, 2).Value = 0.125
$ws.Cells.Item(30, 2).Value = 0.125
$ws.Cells.Item(31, 2).Value = 0.125
$ws.Cells.Item(32, 2).Value = 0.125
$ws.Cells.Item(33, 2).Value = 0.125
$ws.Cells.Item(34, 2).Value = 0.125
$ws.Cells.Item(35, 2).Value = 0.125
$ws.Cells.Item(36, 2).Value = 0.125
$ws.Cells.Item(37, 2).Value = 0.125
$ws.Cells.Item(38, 2).Value = 0.125
$ws.Cells.Item(39, 2).Value = 0.125
$ws.Cells.Item(40, 2).Value = 0.125
$ws.Cells.Item(41, 2).Value = 0.125
$ws.Cells.Item(42, 2).Value = 0.125
$ws.Cells.Item(43, 2).Value = 0.125
$ws.Cells.Item(44, 2).Value = 0.125
$ws.Cells.Item(45, 2).Value = 0.125
$ws.Cells.Item(46, 2).Value = 0.125
$ws.Cells.Item(47, 2).Value = 0.125
$ws.Cells.Item(48, 2).Value = 0.125
$ws.Cells.Item(49, 2).Value = 0.125
$ws.Cells.Item(50, 2).Value = 0.125
$ws.Cells.Item(51, 2).Value = 0.125
$ws.Cells.Item(52, 2).Value = 0.125
$ws.Cells.Item(53, 2).Value = 0.125
$ws.Cells.Item(54, 2).Value = 0.125
$ws.Cells.Item(55, 2).Value = 0.125
$ws.Cells.Item(56, 2).Value = 0.125
$ws.Cells.Item(57, 2).Value = 0.125
$ws.Cells.Item(58, 2).Value = 0.125
$ws.Cells.Item(59, 2).Value = 0.125
$ws.Cells.Item(60, 2).Value = 0.125
$ws.Cells.Item(61, 2).Value = 0.125
$ws.Cells.Item(62, 2).Value = 0.125
$ws.Cells.Item(63, 2).Value = 0.125
$ws.Cells.Item(64, 2).Value = 0.125
$ws.Cells.Item(65, 2).Value = 0.125
$ws.Cells.Item(66, 2).Value = 0.125
$ws.Cells.Item(67, 2).Value = 0.125
$ws.Cells.Item(68, 2).Value = 0.125
$ws.Cells.Item(69, 2).Value = 0.125
$ws.Cells.Item(70, 2).Value = 0.125
$ws.Cells.Item(71, 2).Value = 0.125
$ws.Cells.Item(72, 2).Value = 0.125
$ws.Cells.Item(73, 2).Value = 0.125
$ws.Cells.Item(74, 2).Value = 0.125
$ws.Cells.Item(75, 2).Value = 0.125
$ws.Cells.Item(76, 2).Value = 0.125
$ws.Cells.Item(77, 2).Value = 0.125
$ws.Cells.Item(78, 2).Value = 0.125
$ws.Cells.Item(79, 2).Value = 0.125
$ws.Cells.Item(80, 2).Value = 0.125
$ws.Cells.Item(81, 2).Value = 0.125
$ws.Cells.Item(82, 2).Value = 0.125
$ws.Cells.Item(83, 2).Value = 0.125
$ws.Cells.Item(84, 2).Value = 0.125
$ws.Cells.Item(85, 2).Value = 0.125
$ws.Cells.Item(86, 2).Value = 0.125
$ws.Cells.Item(87, 2).Value = 0.125
$ws.Cells.Item(88, 2).Value = 0.125
$ws.Cells.Item(89, 2).Value = 0.125
$ws.Cells.Item(90, 2).Value = 0.125
$ws.Cells.Item(91, 2).Value = 0.125
$ws.Cells.Item(92, 2).Value = 0.125
$ws.Cells.Item(93, 2).Value = 0.125
$ws.Cells.Item(94, 2).Value = 0.125
$ws.Cells.Item(95, 2).Value = 0.125
$ws.Cells.Item(96, 2).Value = 0.125
$ws.Cells.Item(97, 2).Value = 0.125
$ws.Cells.Item(98, 2).Value = 0.125
$ws.Cells.Item(99, 2).Value = 0.125
$ws.Cells.Item(100, 2).Value = 0.125
$ws.Cells.Item(101, 2).Value = 0.125
$ws.Cells.Item(102, 2).Value = 0.125
$ws.Cells.Item(103, 2).Value = 0.140625
$ws.Cells.Item(104, 2).Value = 0.078125
$ws.Cells.Item(105, 2).Value = 0.125
$ws.Cells.Item(106, 2).Value = 0.125
$ws.Cells.Item(108, 2).Value = 0.078125
$ws.Cells.Item(110, 2).Value = 0.09375
$ws.Cells.Item(111, 2).Value = 0.078125
$ws.Cells.Item(112, 2).Value = 0.078125
$ws.Cells.Item(113, 2).Value = 0.171875
$ws.Cells.Item(114, 2).Value = 0.078125
$ws.Cells.Item(115, 2).Value = 0.0625
$ws.Cells.Item(116, 2).Value = 0.09375
$ws.Cells.Item(118, 2).Value = 0.1475409836065574
